# Update countries & provincias Spain
# Refresh the "paises" COVID dashboard snapshot: bump the "Datos actualizados"
# timestamp, update several countries' case counters, and re-rank two pairs
# of countries (Marruecos now overtakes Croacia; Albania now overtakes Niger)
# whose row order is driven by descending "Casos totales" (column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 11:22"

# --- Row 13 -------------------------------------------------------------
$ws.Range("B13").Value = 26667
$ws.Range("C13").Value = 1684
$ws.Range("D13").Value = 5568
$ws.Range("E13").Value = 18080
$ws.Range("F13").Value = 1278
$ws.Range("G13").Value = 496
$ws.Range("H13").Value = 3019

# --- Row 19 -------------------------------------------------------------
$ws.Range("B19").Value = 13377
$ws.Range("C19").Value = 133
$ws.Range("E19").Value = 6994
$ws.Range("F19").Value = 261

# --- Row 27 -------------------------------------------------------------
$ws.Range("D27").Value = 3141
$ws.Range("E27").Value = 3009
$ws.Range("F27").Value = 74

# --- Rows 60/61: Marruecos overtakes Croacia -----------------------------
$ws.Range("A60").Value = "Marruecos"
$ws.Range("B60").Value = 1431
$ws.Range("C60").Value = 57
$ws.Range("D60").Value = 114
$ws.Range("E60").Value = 1212
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 8
$ws.Range("H60").Value = 105

$ws.Range("A61").Value = "Croacia"
$ws.Range("B61").Value = 1407
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 219
$ws.Range("E61").Value = 1168
$ws.Range("F61").Value = 34
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 20

# --- Row 76 -------------------------------------------------------------
$ws.Range("D76").Value = 63
$ws.Range("E76").Value = 730

# --- Rows 95/96: Albania overtakes Niger ---------------------------------
$ws.Range("A95").Value = "Albania"
$ws.Range("B95").Value = 416
$ws.Range("C95").Value = 7
$ws.Range("D95").Value = 182
$ws.Range("E95").Value = 211
$ws.Range("F95").Value = 6
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 23

$ws.Range("A96").Value = "Niger"
$ws.Range("B96").Value = 410
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 40
$ws.Range("E96").Value = 359
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 11
